$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 575.86664
$ws.Range("I2").Value = 361.2857
$ws.Range("K2").Value = 361.2857
$ws.Range("M2").Value = -248.2857
$ws.Range("H15").Value = 1066.9803
$ws.Range("I15").Value = 1066.9803
$ws.Range("K15").Value = 3200.9409
$ws.Range("M15").Value = -3031.9409
$ws.Range("H18").Value = 9366.429
$ws.Range("I18").Value = 9141.75
$ws.Range("K18").Value = 9141.75
$ws.Range("M18").Value = -8857.75
$ws.Range("H41").Value = 19608572
$ws.Range("I41").Value = 518.4545000000001
$ws.Range("K41").Value = 518.4545000000001
$ws.Range("M41").Value = -78.45450000000005
$ws.Range("H43").Value = 4481.231
$ws.Range("I43").Value = 1793.8572
$ws.Range("K43").Value = 1793.8572
$ws.Range("M43").Value = -1724.8572
$ws.Range("H58").Value = 1048.5385
$ws.Range("I58").Value = 202
$ws.Range("J58").Value = 2036.1666
$ws.Range("K58").Value = 606
$ws.Range("L58").Value = 6108.4998
$ws.Range("M58").Value = -456
$ws.Range("N58").Value = -6408.4998
$ws.Range("H69").Value = 7457.9414
$ws.Range("J69").Value = 7984.643
$ws.Range("L69").Value = 23953.929
$ws.Range("N69").Value = -25701.929
$ws.Range("H72").Value = 7457.9414
$ws.Range("J72").Value = 7984.643
$ws.Range("L72").Value = 71861.787
$ws.Range("N72").Value = -80597.787
$ws.Range("H80").Value = 6049.875
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 8419.799999999999
$ws.Range("K80").Value = 6300
$ws.Range("L80").Value = 25259.4
$ws.Range("M80").Value = -5302
$ws.Range("N80").Value = -27255.4
$ws.Range("H83").Value = 6049.875
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 8419.799999999999
$ws.Range("K83").Value = 18900
$ws.Range("L83").Value = 75778.2
$ws.Range("M83").Value = -13908
$ws.Range("N83").Value = -85762.2
$ws.Range("H86").Value = 8764.695
$ws.Range("I86").Value = 7900.125
$ws.Range("K86").Value = 7900.125
$ws.Range("M86").Value = -6777.125
$ws.Range("H89").Value = 8764.695
$ws.Range("I89").Value = 7900.125
$ws.Range("K89").Value = 39500.625
$ws.Range("M89").Value = -33884.625
$ws.Range("H98").Value = 3345.7273
$ws.Range("I98").Value = 3345.7273
$ws.Range("K98").Value = 3345.7273
$ws.Range("M98").Value = -1847.7273
$ws.Range("H122").Value = 3345.7273
$ws.Range("I122").Value = 3345.7273
$ws.Range("K122").Value = 10037.1819
$ws.Range("M122").Value = -7587.1819
$ws.Range("H131").Value = 8098.5
$ws.Range("I131").Value = 4666.3335
$ws.Range("K131").Value = 13999.0005
$ws.Range("M131").Value = -8959.000499999998
$ws.Range("H137").Value = 70551.96000000001
$ws.Range("I137").Value = 95645
$ws.Range("K137").Value = 286935
$ws.Range("M137").Value = -284385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8885.433000000001
$ws.Range("I32").Value = 5000.3403
$ws.Range("J32").Value = 18015.4
$ws.Range("K32").Value = 5000.3403
$ws.Range("L32").Value = 18015.4
$ws.Range("M32").Value = -4713.3403
$ws.Range("N32").Value = -18589.4
$ws.Range("H45").Value = 4362375
$ws.Range("I45").Value = 6540574.5
$ws.Range("K45").Value = 6540574.5
$ws.Range("M45").Value = -6540197.5
$ws.Range("H122").Value = 7598146.5
$ws.Range("I122").Value = 10583809
$ws.Range("J122").Value = 1898245
$ws.Range("K122").Value = 31751427
$ws.Range("L122").Value = 5694735
$ws.Range("M122").Value = -31748977
$ws.Range("N122").Value = -5699635

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 55560124
$ws.Range("I20").Value = 166667680
$ws.Range("J20").Value = 6350
$ws.Range("K20").Value = 166667680
$ws.Range("L20").Value = 6350
$ws.Range("M20").Value = -166667433
$ws.Range("N20").Value = -6844
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3245.4666
$ws.Range("I99").Value = 2541.1428
$ws.Range("J99").Value = 3861.75
$ws.Range("K99").Value = 2541.1428
$ws.Range("L99").Value = 3861.75
$ws.Range("M99").Value = -1043.1428
$ws.Range("N99").Value = -6857.75
$ws.Range("H107").Value = 34485260
$ws.Range("I107").Value = 2056.647
$ws.Range("K107").Value = 2056.647
$ws.Range("M107").Value = -136.6469999999999
$ws.Range("H122").Value = 1815.5927
$ws.Range("I122").Value = 1638.5416
$ws.Range("K122").Value = 4915.6248
$ws.Range("M122").Value = -2465.6248
$ws.Range("H126").Value = 3245.4666
$ws.Range("I126").Value = 2541.1428
$ws.Range("J126").Value = 3861.75
$ws.Range("K126").Value = 7623.428400000001
$ws.Range("L126").Value = 11585.25
$ws.Range("M126").Value = -5153.428400000001
$ws.Range("N126").Value = -16525.25
$ws.Range("H132").Value = 74481.30499999999
$ws.Range("I132").Value = 49544.57
$ws.Range("K132").Value = 148633.71
$ws.Range("M132").Value = -146103.71
$ws.Range("H141").Value = 504163
$ws.Range("J141").Value = 504163
$ws.Range("L141").Value = 504163
$ws.Range("N141").Value = -514523

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 72973.36
$ws.Range("I5").Value = 1203.5
$ws.Range("J5").Value = 168666.5
$ws.Range("K5").Value = 3610.5
$ws.Range("L5").Value = 505999.5
$ws.Range("M5").Value = -3498.5
$ws.Range("N5").Value = -506223.5
$ws.Range("H6").Value = 102
$ws.Range("I6").Value = 102
$ws.Range("K6").Value = 306
$ws.Range("M6").Value = -193
$ws.Range("H135").Value = 72973.36
$ws.Range("I135").Value = 1203.5
$ws.Range("J135").Value = 168666.5
$ws.Range("K135").Value = 10831.5
$ws.Range("L135").Value = 1517998.5
$ws.Range("M135").Value = -8296.5
$ws.Range("N135").Value = -1523068.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 144.5
$ws.Range("K2").Value = 144.5
$ws.Range("M2").Value = -31.5
$ws.Range("H70").Value = 28589858
$ws.Range("I70").Value = 33353668
$ws.Range("K70").Value = 33353668
$ws.Range("M70").Value = -33353398
$ws.Range("H73").Value = 28589858
$ws.Range("I73").Value = 33353668
$ws.Range("K73").Value = 33353668
$ws.Range("M73").Value = -33352732
$ws.Range("H80").Value = 65551132
$ws.Range("I80").Value = 87400540
$ws.Range("J80").Value = 2899.5
$ws.Range("K80").Value = 87400540
$ws.Range("L80").Value = 2899.5
$ws.Range("M80").Value = -87399542
$ws.Range("N80").Value = -4895.5
$ws.Range("H83").Value = 65551132
$ws.Range("I83").Value = 87400540
$ws.Range("J83").Value = 2899.5
$ws.Range("K83").Value = 437002700
$ws.Range("L83").Value = 14497.5
$ws.Range("M83").Value = -436997708
$ws.Range("N83").Value = -24481.5
$ws.Range("H102").Value = 3487852.2
$ws.Range("I102").Value = 5557191
$ws.Range("K102").Value = 5557191
$ws.Range("M102").Value = -5555569
$ws.Range("H122").Value = 279730.56
$ws.Range("I122").Value = 425142.1
$ws.Range("K122").Value = 1275426.3
$ws.Range("M122").Value = -1272976.3
$ws.Range("H126").Value = 3070180.8
$ws.Range("I126").Value = 2395246.8
$ws.Range("J126").Value = 3627735
$ws.Range("K126").Value = 7185740.399999999
$ws.Range("L126").Value = 10883205
$ws.Range("M126").Value = -7183270.399999999
$ws.Range("N126").Value = -10888145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6356.5
$ws.Range("I7").Value = 4555.125
$ws.Range("K7").Value = 4555.125
$ws.Range("M7").Value = -4443.125
$ws.Range("H126").Value = 6356.5
$ws.Range("I126").Value = 4555.125
$ws.Range("K126").Value = 13665.375
$ws.Range("M126").Value = -11195.375
$ws.Range("H132").Value = 3460.0444
$ws.Range("I132").Value = 3081.4482
$ws.Range("J132").Value = 4146.25
$ws.Range("K132").Value = 9244.3446
$ws.Range("L132").Value = 12438.75
$ws.Range("M132").Value = -6714.3446
$ws.Range("N132").Value = -17498.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18522426
$ws.Range("I81").Value = 23813992
$ws.Range("J81").Value = 1950
$ws.Range("K81").Value = 47627984
$ws.Range("L81").Value = 3900
$ws.Range("M81").Value = -47626923
$ws.Range("N81").Value = -6022
$ws.Range("H84").Value = 18522426
$ws.Range("I84").Value = 23813992
$ws.Range("J84").Value = 1950
$ws.Range("K84").Value = 238139920
$ws.Range("L84").Value = 19500
$ws.Range("M84").Value = -238134616
$ws.Range("N84").Value = -30108
$ws.Range("H122").Value = 2409.889
$ws.Range("I122").Value = 1711.5385
$ws.Range("J122").Value = 3365.5264
$ws.Range("K122").Value = 5134.6155
$ws.Range("L122").Value = 10096.5792
$ws.Range("M122").Value = -2684.6155
$ws.Range("N122").Value = -14996.5792
$ws.Range("H126").Value = 1871.0358
$ws.Range("I126").Value = 1978.5294
$ws.Range("K126").Value = 5935.5882
$ws.Range("M126").Value = -3465.5882
